$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (180 and 181) to the feed log data.
$ws.Cells.Item(180, 1).Value = 179
$ws.Cells.Item(180, 2).Value = 1
$ws.Cells.Item(180, 3).Value = "2024-06-18 18:17:31"
$ws.Cells.Item(180, 4).Value = 200
$ws.Cells.Item(180, 5).Value = 19

$ws.Cells.Item(181, 1).Value = 180
$ws.Cells.Item(181, 2).Value = 2
$ws.Cells.Item(181, 3).Value = "2024-06-18 18:17:32"
$ws.Cells.Item(181, 4).Value = 200
$ws.Cells.Item(181, 5).Value = 1
